$p = $ppt.ActivePresentation
$s = $p.Slides.Item(13)
$tr = $s.Shapes.Item(2).TextFrame.TextRange

# Fix 1: merge "of your " + "recital" runs into one run "of your recital"
$para3 = $tr.Paragraphs(3,1)
$chars = $para3.Characters(22, 15)
$chars.Text = "of your recital"

# Fix 2: merge "Various "+"level of Swara Reach – Show "+"table" runs into one run,
# and remove the trailing endParaRPr by re-creating the paragraph at the end.
$para5 = $tr.Paragraphs(5,1)
$para5.InsertAfter("`rVarious level of Swara Reach – Show table")
$tr.Paragraphs(5,1).Delete()

"done"
